# Update "想去人数" (interest count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 150
$ws1.Range("F18").Value = 5099
$ws1.Range("F21").Value = 117
$ws1.Range("F22").Value = 2274
$ws1.Range("F24").Value = 30
$ws1.Range("F25").Value = 2122

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 150
$ws4.Range("F18").Value = 5099
$ws4.Range("F23").Value = 117
$ws4.Range("F24").Value = 2274
$ws4.Range("F26").Value = 30
$ws4.Range("F27").Value = 2122
